# Generate Report for Handoff
# Adds a new row (row 9) to each of the three worksheets (Overview, zh-cn, de-de)
# describing the newly handed-off file
# "d195734e-17e9-4871-92ec-508ae62431ab.md".

$wb = $excel.ActiveWorkbook

$hyperlinkUnderline = 2            # xlUnderlineStyleSingle
$hyperlinkColor = 15570276         # BGR encoding of RGB(0x64,0x95,0xED) used by the existing HyperLink style

function Set-HyperlinkLook($range) {
    $range.Font.Underline = $hyperlinkUnderline
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A9").Value = "d195734e-17e9-4871-92ec-508ae62431ab.md"
Set-HyperlinkLook $wsOverview.Range("A9")

$wsOverview.Range("B9").Value = "Ready for handoff"
$wsOverview.Range("C9").Value = "Ready for handoff"
$wsOverview.Range("D9").Value = "2016-31-17 12:31:13"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f2c9a61c3a1e4d9fa0c1a5b1e5b7c0a4b9d4e2f1/e2e/d195734e-17e9-4871-92ec-508ae62431ab.md",
    "",
    "",
    "d195734e-17e9-4871-92ec-508ae62431ab.md"
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A9").Value = "d195734e-17e9-4871-92ec-508ae62431ab.md"
Set-HyperlinkLook $wsZhCn.Range("A9")

$wsZhCn.Range("B9").Value = ".md"
Set-HyperlinkLook $wsZhCn.Range("B9")

$wsZhCn.Range("C9").Value = "Ready for handoff"

$wsZhCn.Range("D9").Value = "d195734e-17e9-4871-92ec-508ae62431ab.f1f0ba68e2eb53ee3fe8b3b7ed35dc5a369f67a7.zh-cn.xlf"
Set-HyperlinkLook $wsZhCn.Range("D9")

$wsZhCn.Range("E9").Value = "2016-03-17 12:31:07"
$wsZhCn.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Range("H9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I9").Value = "Include"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f2c9a61c3a1e4d9fa0c1a5b1e5b7c0a4b9d4e2f1/e2e/d195734e-17e9-4871-92ec-508ae62431ab.md",
    "",
    "",
    "d195734e-17e9-4871-92ec-508ae62431ab.md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("B9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f2c9a61c3a1e4d9fa0c1a5b1e5b7c0a4b9d4e2f1/e2e/d195734e-17e9-4871-92ec-508ae62431ab.md",
    "",
    "",
    ".md"
) | Out-Null

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fe9ba81990d528bc59d5c0060634a3cc44edae04/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d195734e-17e9-4871-92ec-508ae62431ab.f1f0ba68e2eb53ee3fe8b3b7ed35dc5a369f67a7.zh-cn.xlf",
    "",
    "",
    "d195734e-17e9-4871-92ec-508ae62431ab.f1f0ba68e2eb53ee3fe8b3b7ed35dc5a369f67a7.zh-cn.xlf"
) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A9").Value = "d195734e-17e9-4871-92ec-508ae62431ab.md"
Set-HyperlinkLook $wsDeDe.Range("A9")

$wsDeDe.Range("B9").Value = ".md"
Set-HyperlinkLook $wsDeDe.Range("B9")

$wsDeDe.Range("C9").Value = "Ready for handoff"

$wsDeDe.Range("D9").Value = "d195734e-17e9-4871-92ec-508ae62431ab.f1f0ba68e2eb53ee3fe8b3b7ed35dc5a369f67a7.de-de.xlf"
Set-HyperlinkLook $wsDeDe.Range("D9")

$wsDeDe.Range("E9").Value = "2016-03-17 12:31:13"
$wsDeDe.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Range("H9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I9").Value = "Include"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f2c9a61c3a1e4d9fa0c1a5b1e5b7c0a4b9d4e2f1/e2e/d195734e-17e9-4871-92ec-508ae62431ab.md",
    "",
    "",
    "d195734e-17e9-4871-92ec-508ae62431ab.md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("B9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/f2c9a61c3a1e4d9fa0c1a5b1e5b7c0a4b9d4e2f1/e2e/d195734e-17e9-4871-92ec-508ae62431ab.md",
    "",
    "",
    ".md"
) | Out-Null

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b2e3abcb3c700dddb9f29c46adcbd2eec8ed51e6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d195734e-17e9-4871-92ec-508ae62431ab.f1f0ba68e2eb53ee3fe8b3b7ed35dc5a369f67a7.de-de.xlf",
    "",
    "",
    "d195734e-17e9-4871-92ec-508ae62431ab.f1f0ba68e2eb53ee3fe8b3b7ed35dc5a369f67a7.de-de.xlf"
) | Out-Null

Write-Output "Report row added to Overview, zh-cn and de-de sheets"
